$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Plans")

# Switch the two test-plan rows from the single shared "Reactivate" name to
# distinct names, and point the CSV Data Set file name (column AC) at the
# new relative path instead of the old machine-local absolute path.
$ws.Range("AC2").Value = "src/test/resources/data/create_jmx_files/csv_data/data_driven_create_jmx_file.csv"
$ws.Range("AC3").Value = "src/test/resources/data/create_jmx_files/csv_data/data_driven_create_jmx_file.csv"

$ws.Range("A2").Value = "Test Plan One"
$ws.Range("A3").Value = "Test Plan Two"

# Update the sheet view: drop the frozen/scrolled topLeftCell and move the
# active selection back to A3.
$ws.Application.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("A3").Select()
